$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p002_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p002r_1</id>", 2)
$d.Content.Find.Execute("<id>p002r_2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p002r_2</id>", 2)
